# Clean up footnote markers (e.g. " [5]") and embedded line breaks in the
# Vaccine/BrandName label strings across every sheet of the workbook.
#
# The source data used a trailing "[N]" citation marker on many vaccine
# category names, and some long labels (and a couple of brand names) were
# hard-wrapped with an embedded newline. This removes the "[N]" markers and
# collapses the embedded newlines into a single space, matching the cleaned
# -up label text used everywhere else in the sheet.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rowCount = $used.Rows.Count
    $colCount = $used.Columns.Count
    $rowStart = $used.Row
    $colStart = $used.Column

    for ($r = $rowStart; $r -lt ($rowStart + $rowCount); $r++) {
        for ($c = $colStart; $c -lt ($colStart + $colCount); $c++) {
            $cell = $ws.Cells.Item($r, $c)
            $val = $cell.Value2

            if ($val -is [string]) {
                $newVal = $val -replace "`n", ' '
                $newVal = $newVal -replace '\[\d+\]', ''

                if ($newVal -ne $val) {
                    $cell.Value2 = $newVal
                }
            }
        }
    }
}
